$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '71.287.31'
Set-TextValue $ws.Range('E2') '  -2.48%  '
Set-TextValue $ws.Range('D3') '3.871.84'
Set-TextValue $ws.Range('E3') '  -3.12%  '
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('E4') '  -0.08%  '
Set-TextValue $ws.Range('D5') '601.68'
Set-TextValue $ws.Range('E5') '  +0.65%  '
Set-TextValue $ws.Range('D6') '170.19'
Set-TextValue $ws.Range('E6') '  +3.34%  '
Set-TextValue $ws.Range('D7') '0.669'
Set-TextValue $ws.Range('E7') '  -2.29%  '
Set-TextValue $ws.Range('E8') '  +0.11%  '
Set-TextValue $ws.Range('D9') '0.748'
Set-TextValue $ws.Range('E9') '  -0.62%  '
Set-TextValue $ws.Range('D10') '0.178'
Set-TextValue $ws.Range('E10') '  +5.03%  '
Set-TextValue $ws.Range('D11') '53.62'
Set-TextValue $ws.Range('E11') '  -2.67%  '
Set-TextValue $ws.Range('D12') '0.0000323'
Set-TextValue $ws.Range('E12') '  +0.94%  '
Set-TextValue $ws.Range('D13') '11.50'
Set-TextValue $ws.Range('E13') '  +4.39%  '
Set-TextValue $ws.Range('D14') '4.470.94'
Set-TextValue $ws.Range('E14') '  -3.55%  '
Set-TextValue $ws.Range('D15') '21.19'
Set-TextValue $ws.Range('E15') '  +3.26%  '
Set-TextValue $ws.Range('D16') '3.856.67'
Set-TextValue $ws.Range('E16') '  -3.61%  '
Set-TextValue $ws.Range('D17') '13.95'
Set-TextValue $ws.Range('E17') '  -1.33%  '
Set-TextValue $ws.Range('E18') '  -4.21%  '
Set-TextValue $ws.Range('E19') '  -2.38%  '
Set-TextValue $ws.Range('D20') '71.023.04'
Set-TextValue $ws.Range('E20') '  -2.48%  '
Set-TextValue $ws.Range('D21') '440.74'
Set-TextValue $ws.Range('E21') '  +0.46%  '
Set-TextValue $ws.Range('D22') '4.79'
Set-TextValue $ws.Range('E22') '  +0.99%  '
Set-TextValue $ws.Range('D23') '94.81'
Set-TextValue $ws.Range('E23') '  -2.02%  '
Set-TextValue $ws.Range('E24') '  -4.66%  '
Set-TextValue $ws.Range('D25') '13.87'
Set-TextValue $ws.Range('E25') '  -3.56%  '
Set-TextValue $ws.Range('D26') '11.65'
Set-TextValue $ws.Range('E26') '  +2.36%  '
Set-TextValue $ws.Range('D27') '4.01'
Set-TextValue $ws.Range('E27') '  -7.61%  '
Set-TextValue $ws.Range('D28') '5.98'
Set-TextValue $ws.Range('E28') '  +0.18%  '
Set-TextValue $ws.Range('D29') '10.55'
Set-TextValue $ws.Range('E29') '  +1.39%  '
Set-TextValue $ws.Range('D30') '8.64'
Set-TextValue $ws.Range('E30') '  +8.63%  '
Set-TextValue $ws.Range('D31') '35.14'
Set-TextValue $ws.Range('E31') '  -3.61%  '
Set-TextValue $ws.Range('D32') '13.56'
Set-TextValue $ws.Range('E32') '  -2.54%  '
Set-TextValue $ws.Range('D33') '48.24'
Set-TextValue $ws.Range('D34') '0.126'
Set-TextValue $ws.Range('E34') '  -4.03%  '
Set-TextValue $ws.Range('B35') 'PEPE'
Set-TextValue $ws.Range('C35') 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D35') '0.0₃0989'
Set-TextValue $ws.Range('E35') '  +9.28%  '
Set-TextValue $ws.Range('B36') 'OKB'
Set-TextValue $ws.Range('C36') 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range('D36') '68.99'
Set-TextValue $ws.Range('E36') '  -2.43%  '
Set-TextValue $ws.Range('D37') '640.73'
Set-TextValue $ws.Range('E37') '  -4.54%  '
Set-TextValue $ws.Range('D38') '0.440'
Set-TextValue $ws.Range('E38') '  +0.35%  '
Set-TextValue $ws.Range('D39') '0.147'
Set-TextValue $ws.Range('E39') '  +0.31%  '
Set-TextValue $ws.Range('E40') '  +0.40%  '
Set-TextValue $ws.Range('D41') '0.999'
Set-TextValue $ws.Range('E41') '  -0.36%  '
Set-TextValue $ws.Range('D42') '3.26'
Set-TextValue $ws.Range('E42') '  -3.05%  '
Set-TextValue $ws.Range('D43') '2.86'
Set-TextValue $ws.Range('E43') '  +8.84%  '
Set-TextValue $ws.Range('D44') '0.0472'
Set-TextValue $ws.Range('E44') '  -3.84%  '
Set-TextValue $ws.Range('D45') '3.10'
Set-TextValue $ws.Range('E45') '  +15.07%  '
Set-TextValue $ws.Range('D46') '10.28'
Set-TextValue $ws.Range('E46') '  -4.06%  '
Set-TextValue $ws.Range('B47') 'Stellar'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D47') '0.144'
Set-TextValue $ws.Range('E47') '  -3.63%  '
Set-TextValue $ws.Range('B48') 'WEMIXToken'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D48') '2.91'
Set-TextValue $ws.Range('E48') '  -13.40%  '
Set-TextValue $ws.Range('D49') '2.966.22'
Set-TextValue $ws.Range('E49') '  +1.68%  '
Set-TextValue $ws.Range('D50') '3.29'
Set-TextValue $ws.Range('E50') '  -3.36%  '
Set-TextValue $ws.Range('D51') '0.000279'
Set-TextValue $ws.Range('E51') '  +3.91%  '
